$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
# Row 125
$ws.Range("H125").Value = 19232712
$ws.Range("I125").Value = 20000620
$ws.Range("J125").Value = 35000
$ws.Range("K125").Value = 180005580
$ws.Range("L125").Value = 315000
$ws.Range("M125").Value = -180003120
$ws.Range("N125").Value = -319920
# Row 138
$ws.Range("H138").Value = 4177.1
$ws.Range("I138").Value = 1357.8572
$ws.Range("J138").Value = 4636.0464
$ws.Range("K138").Value = 4073.5716
$ws.Range("L138").Value = 13908.1392
$ws.Range("M138").Value = 1066.4284
$ws.Range("N138").Value = -24188.1392

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 17470.229
$ws.Range("I32").Value = 17453.273
$ws.Range("K32").Value = 17453.273
$ws.Range("M32").Value = -17166.273
# Row 45
$ws.Range("H45").Value = 5691.577
$ws.Range("I45").Value = 7163.353
$ws.Range("J45").Value = 2911.5557
$ws.Range("K45").Value = 7163.353
$ws.Range("L45").Value = 2911.5557
$ws.Range("M45").Value = -6786.353
$ws.Range("N45").Value = -3665.5557
# Row 110
$ws.Range("H110").Value = 657.5833
$ws.Range("I110").Value = 626.4545000000001
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 626.4545000000001
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = 1418.5455
$ws.Range("N110").Value = -5090

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2388.4167
$ws.Range("I86").Value = 2301.2
$ws.Range("J86").Value = 2533.7778
$ws.Range("K86").Value = 2301.2
$ws.Range("L86").Value = 2533.7778
$ws.Range("M86").Value = -1178.2
$ws.Range("N86").Value = -4779.7778
# Row 89
$ws.Range("H89").Value = 2388.4167
$ws.Range("I89").Value = 2301.2
$ws.Range("J89").Value = 2533.7778
$ws.Range("K89").Value = 11506
$ws.Range("L89").Value = 12668.889
$ws.Range("M89").Value = -5890
$ws.Range("N89").Value = -23900.889
# Row 107
$ws.Range("H107").Value = 1328
$ws.Range("I107").Value = 1262.2
$ws.Range("J107").Value = 1437.6666
$ws.Range("K107").Value = 1262.2
$ws.Range("L107").Value = 1437.6666
$ws.Range("M107").Value = 657.8
$ws.Range("N107").Value = -5277.6666

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
# Row 41
$ws.Range("H41").Value = 2116.3333
$ws.Range("I41").Value = 2116.3333
$ws.Range("K41").Value = 2116.3333
$ws.Range("M41").Value = -1688.3333
# Row 58
$ws.Range("H58").Value = 7725.1665
$ws.Range("I58").Value = 2120.7646
$ws.Range("J58").Value = 103000
$ws.Range("K58").Value = 2120.7646
$ws.Range("L58").Value = 103000
$ws.Range("M58").Value = -1917.7646
$ws.Range("N58").Value = -103406
# Row 59
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
# Row 136
$ws.Range("H136").Value = 7725.1665
$ws.Range("I136").Value = 2120.7646
$ws.Range("J136").Value = 103000
$ws.Range("K136").Value = 6362.293799999999
$ws.Range("L136").Value = 309000
$ws.Range("M136").Value = -3812.293799999999
$ws.Range("N136").Value = -314100

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1066.51
$ws.Range("I68").Value = 610.3939
$ws.Range("J68").Value = 1291.1642
$ws.Range("K68").Value = 1831.1817
$ws.Range("L68").Value = 3873.4926
$ws.Range("M68").Value = -1020.1817
$ws.Range("N68").Value = -5495.4926
# Row 71
$ws.Range("H71").Value = 1066.51
$ws.Range("I71").Value = 610.3939
$ws.Range("J71").Value = 1291.1642
$ws.Range("K71").Value = 5493.5451
$ws.Range("L71").Value = 11620.4778
$ws.Range("M71").Value = -1437.5451
$ws.Range("N71").Value = -19732.4778
# Row 86
$ws.Range("H86").Value = 800
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
# Row 89
$ws.Range("H89").Value = 800
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
# Row 107
$ws.Range("H107").Value = 788.9434
$ws.Range("I107").Value = 254.75
$ws.Range("J107").Value = 1112.697
$ws.Range("K107").Value = 764.25
$ws.Range("L107").Value = 3338.090999999999
$ws.Range("M107").Value = 1155.75
$ws.Range("N107").Value = -7178.090999999999
# Row 118
$ws.Range("H118").Value = 4746.4
$ws.Range("I118").Value = 1150
$ws.Range("J118").Value = 7144
$ws.Range("K118").Value = 3450
$ws.Range("L118").Value = 21432
$ws.Range("M118").Value = -2207
$ws.Range("N118").Value = -23918
# Row 125
$ws.Range("H125").Value = 2681.3
$ws.Range("I125").Value = 920
$ws.Range("J125").Value = 2992.1177
$ws.Range("K125").Value = 2760
$ws.Range("L125").Value = 8976.3531
$ws.Range("M125").Value = 2160
$ws.Range("N125").Value = -18816.3531
# Row 131
$ws.Range("H131").Value = 21156.96
$ws.Range("I131").Value = 72634.28999999999
$ws.Range("J131").Value = 2191.6316
$ws.Range("K131").Value = 217902.87
$ws.Range("L131").Value = 6574.8948
$ws.Range("M131").Value = -212862.87
$ws.Range("N131").Value = -16654.8948

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1462
$ws.Range("I7").Value = 1477.7142
$ws.Range("J7").Value = 1451
$ws.Range("K7").Value = 1477.7142
$ws.Range("L7").Value = 1451
$ws.Range("M7").Value = -1365.7142
$ws.Range("N7").Value = -1675
# Row 122
$ws.Range("H122").Value = 3538.158
$ws.Range("I122").Value = 3816.923
$ws.Range("J122").Value = 2934.1667
$ws.Range("K122").Value = 11450.769
$ws.Range("L122").Value = 8802.500100000001
$ws.Range("M122").Value = -9000.769
$ws.Range("N122").Value = -13702.5001
# Row 126
$ws.Range("H126").Value = 1462
$ws.Range("I126").Value = 1477.7142
$ws.Range("J126").Value = 1451
$ws.Range("K126").Value = 4433.142599999999
$ws.Range("L126").Value = 4353
$ws.Range("M126").Value = -1963.142599999999
$ws.Range("N126").Value = -9293

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 760
$ws.Range("I126").Value = 540
$ws.Range("J126").Value = 1200
$ws.Range("K126").Value = 1620
$ws.Range("L126").Value = 3600
$ws.Range("M126").Value = 850
$ws.Range("N126").Value = -8540
# Row 136
$ws.Range("H136").Value = 1926.2264
$ws.Range("I136").Value = 1987.1875
$ws.Range("J136").Value = 1833.3334
$ws.Range("K136").Value = 5961.5625
$ws.Range("L136").Value = 5500.0002
$ws.Range("M136").Value = -3411.5625
$ws.Range("N136").Value = -10600.0002
